# Update the Metadata sheet:
#  - B7 ("Experimental" row) gets the literal text value "true"
#  - B8 ("Date" row) gets updated to the new timestamp
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Write "true" as a formula that evaluates to the text "true", then paste
# the result back as a value so the cell ends up holding literal text
# (not a Boolean) while keeping its existing style.
$cell = $ws.Range("B7")
$cell.Formula = "=""true"""
$cell.Copy()
$cell.PasteSpecial(-4163)

$ws.Range("B8").Value = "2023-02-16T14:43:10-06:00"
